# Update the "想去人数" (F) and "最低票价" (G) numbers on the "展览" and
# "全部类型" worksheets to reflect the latest scraped counts.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1177
    $ws.Range("F6").Value = 172
    $ws.Range("F10").Value = 5476
    $ws.Range("G10").Value = 70
    $ws.Range("F11").Value = 4880
    $ws.Range("F12").Value = 19
    $ws.Range("F16").Value = 200
}
